$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the "profile" column (currently column U) to before column F,
# shifting the intervening columns (F..T) one place to the right.
$ws.Columns("U:U").Cut()
$ws.Columns("F:F").Insert()

# Give the new (moved) column F its own explicit width
# (target stored width 18.7109375 chars; 17.8 is the closest COM
# ColumnWidth input this engine's pixel-grid rounding reproduces).
$ws.Columns("F").ColumnWidth = 17.8

# Update selection to F9
$ws.Range("F9").Select()
